$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new "Po No" header in column G, matching the style of the other header cells
$ws.Range("G1").Value = "Po No"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active cell selection on the sheet, as saved in the workbook
$ws.Range("G6").Select() | Out-Null

Write-Host "grn import column added"
